# Update column F (dSF) values on Sheet1 to match repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value  = -4
$ws.Range("F5").Value  = -8
$ws.Range("F6").Value  = 10
$ws.Range("F7").Value  = -5
$ws.Range("F9").Value  = -9
$ws.Range("F12").Value = 6
$ws.Range("F14").Value = 3
$ws.Range("F15").Value = -2
$ws.Range("F16").Value = -1
